$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed cryptos snapshot.
# Column D cells are forced to Text format so numeric-looking strings (e.g. "1.002",
# "0.08032") are preserved exactly as text, matching the source inlineStr cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.367.59"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.878.98"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7219"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.15"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08032"
$ws.Range("E8").Value = "  +3.06%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08160"
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.880.46"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.58"
$ws.Range("E13").Value = "  +3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.230"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7109"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.406"
$ws.Range("E16").Value = "  +5.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008478"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.371.24"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.33"
$ws.Range("E19").Value = "  +1.63%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.128.12"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.731"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1610"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.70"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.037"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.280"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.224"
$ws.Range("E32").Value = "  -5.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05344"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.934"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7599"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.176"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.700"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01869"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.264.37"
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.760"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.432"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "112.98"
$ws.Range("E42").Value = "  +3.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9041"
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.03"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("E45").Value = "  +5.86%  "
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.024.83"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.797"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5199"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.473"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4337"
$ws.Range("E51").Value = "  -0.19%  "
